# Fix typo: "zugeen" -> "zugehen" in the verb list (cell A13)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "zugehen"

# Update the selection to match the saved state (user last selected C15)
$ws.Range("C15").Select()
